$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reverse the order of the "Periodo Mora" column (E16:E28), newest period first.
$periods = @(2105, 2104, 2103, 2102, 2101, 2012, 2011, 2010, 2009, 2008, 2007, 2006, 2005)
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = [string]$periods[$i]
}

# The "Valor Mora" values travel with their period: 2105 carries 32000, 2005 carries 40000.
$ws.Cells.Item(16, 6).Value = 32000
$ws.Cells.Item(28, 6).Value = 40000
